# Applies the tracked changes described by the commit diff:
#  1. Insert a new (highlighted) paragraph about AI / data-mining right
#     after the paragraph ending "...hay fever."
#  2. Insert a new (highlighted) paragraph about the project being
#     unlikely to replace GP jobs into the previously-empty paragraph
#     that follows the CSIRO/Roy-Morgan paragraph.
#  3. Remove the stray <w:lastRenderedPageBreak/> that used to sit in
#     front of "HOW WILL THIS AFFECT YOU?" (forcing Word to rebuild that
#     run via Find/Replace drops the stale rendering bookmark).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. New paragraph after "...hich in this case would be hay fever."
# ---------------------------------------------------------------------
$anchor1 = $d.Content
$found1 = $anchor1.Find.Execute(
    "hich in this case would be hay fever.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found1) {
    $anchor1.Collapse(0)
    $anchor1.InsertParagraphAfter()
    $insertPos = $anchor1.Start
    $newPara1 = $d.Range($insertPos, $insertPos)

    $aiText = "This project would be possible by using advanced artificial intelligence (AI) and data mining. Many companies such as Google, Amazon, and Microsoft are working on advanced AI which can be used in many ways and has almost endless possibilities. These AI" + [char]8217 + "s that are being developed alongside data mining or " + [char]8216 + "big data" + [char]8217 + " technology could be used to help the system recognise patterns and medical conditions so that it can provide users with possible conditions."

    # Insert the plain text first (highlighting a freshly-typed Range
    # directly does not stick), then use Find/Replace to paint the
    # highlight onto the text we just wrote, plus a trailing
    # non-highlighted space to match the source paragraph.
    $newPara1.Text = $aiText + " "

    $find1 = $d.Content.Find
    $find1.ClearFormatting()
    $find1.Replacement.ClearFormatting()
    $find1.Replacement.Highlight = $true
    $find1.Execute($aiText, $false, $false, $false, $false, $false, $true, 1, $false, $aiText, 2) | Out-Null
}

# ---------------------------------------------------------------------
# 2. Fill in the empty paragraph after the CSIRO / Roy Morgan paragraph
# ---------------------------------------------------------------------
$anchor2 = $d.Content
$found2 = $anchor2.Find.Execute(
    "without any medical information to back their claims.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found2) {
    $anchor2.Collapse(0)
    $anchor2.MoveEnd(1, 1)
    $anchor2.Collapse(0)
    $targetPos = $anchor2.Start
    $newPara2 = $d.Range($targetPos, $targetPos)

    $unlikelyText = "This project is unlikely to take jobs from health sector or make them redundant as general practitioners would still be a very important part of the health industry. The prediction tool was designed to reduce the amount of people that went to visit their doctors for minor issues, but many people would still need to visit their GPs for serious issues, or ongoing problems that someone has. In addition, the prediction tool would still advise users to visit their GP for information, if necessary."

    $newPara2.Text = $unlikelyText + " "

    $find2 = $d.Content.Find
    $find2.ClearFormatting()
    $find2.Replacement.ClearFormatting()
    $find2.Replacement.Highlight = $true
    $find2.Execute($unlikelyText, $false, $false, $false, $false, $false, $true, 1, $false, $unlikelyText, 2) | Out-Null
}

# ---------------------------------------------------------------------
# 3. Drop the obsolete lastRenderedPageBreak in front of the
#    "HOW WILL THIS AFFECT YOU?" heading by forcing Word to rebuild the
#    run that currently carries it.
# ---------------------------------------------------------------------
$find3 = $d.Content.Find
$find3.ClearFormatting()
$find3.Replacement.ClearFormatting()
$headingText = "HOW WILL THIS AFFECT YOU?"
$find3.Execute($headingText, $false, $false, $false, $false, $false, $true, 1, $false, $headingText, 2) | Out-Null

Write-Host "done"
